# Update countries & provincias Spain
# - Swap El Salvador / Madagascar order in the country list (rows 130/131)
#   and refresh their stats.
# - Refresh the "last updated" timestamp in A1.
# - Refresh numeric COVID stats for a handful of other countries (rows
#   10,14,17,18,21,24,25,27,48,68,99).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 13:22"

# --- Row 10 (Iran) ------------------------------------------------------
$ws.Range("B10").Value = 64586
$ws.Range("C10").Value = 1997
$ws.Range("E10").Value = 33554
$ws.Range("G10").Value = 121
$ws.Range("H10").Value = 3993

# --- Row 14 (Suiza) ------------------------------------------------------
$ws.Range("B14").Value = 22789
$ws.Range("C14").Value = 536
$ws.Range("E14").Value = 13239

# --- Row 17 (Brasil) ------------------------------------------------------
$ws.Range("B17").Value = 14072
$ws.Range("C17").Value = 38
$ws.Range("E17").Value = 13254
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 691

# --- Row 18 (Austria) ------------------------------------------------------
$ws.Range("B18").Value = 12782
$ws.Range("C18").Value = 143
$ws.Range("E18").Value = 7997

# --- Row 21 (Israel) ------------------------------------------------------
$ws.Range("E21").Value = 8531
$ws.Range("G21").Value = 7
$ws.Range("H21").Value = 72

# --- Row 24 (Noruega) ------------------------------------------------------
$ws.Range("E24").Value = 5961
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 93

# --- Row 25 (Australia) ------------------------------------------------------
$ws.Range("B25").Value = 6013
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 2813
$ws.Range("E25").Value = 3150
$ws.Range("F25").Value = 87

# --- Row 27 (Dinamarca) ------------------------------------------------------
$ws.Range("D27").Value = 1621
$ws.Range("E27").Value = 3547
$ws.Range("G27").Value = 15
$ws.Range("H27").Value = 218

# --- Row 48 (Catar) ------------------------------------------------------
$ws.Range("B48").Value = 2210
$ws.Range("C48").Value = 153
$ws.Range("D48").Value = 178
$ws.Range("E48").Value = 2026

# --- Row 68 (Lituania) ------------------------------------------------------
$ws.Range("F68").Value = 21

# --- Row 99 (Malta) ------------------------------------------------------
$ws.Range("B99").Value = 299
$ws.Range("C99").Value = 6
$ws.Range("E99").Value = 294

# --- Rows 130/131: El Salvador & Madagascar swap positions --------------
# Row 130 was "El Salvador", becomes "Madagascar" with fresh stats.
$ws.Range("A130").Value = "Madagascar"
$ws.Range("B130").Value = 93
$ws.Range("C130").Value = 5
$ws.Range("D130").Value = 11
$ws.Range("E130").Value = 82
$ws.Range("F130").Value = 1
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 0

# Row 131 was "Madagascar", becomes "El Salvador" with fresh stats.
$ws.Range("A131").Value = "El Salvador"
$ws.Range("B131").Value = 93
$ws.Range("C131").Value = 15
$ws.Range("D131").Value = 9
$ws.Range("E131").Value = 79
$ws.Range("F131").Value = 2
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 5
